$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update raw input values (rebalancing) ---
$ws.Range("B2").Value = 2150
$ws.Range("C2").Value = 345

$ws.Range("B3").Value = 27000
$ws.Range("C3").Value = 295

$ws.Range("B4").Value = 89000

$ws.Range("B5").Value = 8.5

# --- Add new formula cell B8 ---
$ws.Range("B8").Formula = "=250/335*18"

# --- Update the active selection to D2 ---
[void]$ws.Range("D2").Select()
